$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these cells to Text format so numeric-looking values (e.g. "7.80")
# are stored verbatim as strings, matching the source data which keeps
# trailing zeros / exact formatting instead of being parsed as numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "64.094.69"
$ws.Range("E2").Value = "  -3.74%  "

# Row 3
$ws.Range("D3").Value = "3.142.94"
$ws.Range("E3").Value = "  -3.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "609.26"
$ws.Range("E5").Value = "  +0.19%  "

# Row 6
$ws.Range("D6").Value = "146.66"
$ws.Range("E6").Value = "  -7.02%  "

# Row 8
$ws.Range("D8").Value = "3.140.72"
$ws.Range("E8").Value = "  -3.44%  "

# Row 9
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -5.01%  "

# Row 10
$ws.Range("E10").Value = "  -7.21%  "

# Row 11
$ws.Range("E11").Value = "  -6.74%  "

# Row 12
$ws.Range("E12").Value = "  -6.02%  "

# Row 13
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  -7.81%  "

# Row 14
$ws.Range("D14").Value = "35.43"
$ws.Range("E14").Value = "  -9.67%  "

# Row 15
$ws.Range("D15").Value = "3.656.74"
$ws.Range("E15").Value = "  -3.45%  "

# Row 16
$ws.Range("D16").Value = "64.123.69"
$ws.Range("E16").Value = "  -3.83%  "

# Row 17
$ws.Range("E17").Value = "  +0.69%  "

# Row 18
$ws.Range("D18").Value = "3.140.22"
$ws.Range("E18").Value = "  -3.57%  "

# Row 19
$ws.Range("E19").Value = "  -7.58%  "

# Row 20
$ws.Range("D20").Value = "477.88"
$ws.Range("E20").Value = "  -5.88%  "

# Row 21
$ws.Range("D21").Value = "14.83"
$ws.Range("E21").Value = "  -3.91%  "

# Row 22
$ws.Range("D22").Value = "0.704"
$ws.Range("E22").Value = "  -6.54%  "

# Row 23
$ws.Range("D23").Value = "7.80"
$ws.Range("E23").Value = "  -3.75%  "

# Row 24
$ws.Range("D24").Value = "13.61"
$ws.Range("E24").Value = "  -7.87%  "

# Row 25
$ws.Range("D25").Value = "83.34"

# Row 26
$ws.Range("E26").Value = "  +0.09%  "

# Row 27
$ws.Range("E27").Value = "  -5.37%  "

# Row 28
$ws.Range("D28").Value = "8.40"
$ws.Range("E28").Value = "  -7.66%  "

# Row 29
$ws.Range("D29").Value = "2.18"

# Row 30
$ws.Range("D30").Value = "6.77"
$ws.Range("E30").Value = "  -1.57%  "

# Row 31
$ws.Range("E31").Value = "  -15.13%  "

# Row 32
$ws.Range("D32").Value = "2.74"
$ws.Range("E32").Value = "  -5.48%  "

# Row 33
$ws.Range("E33").Value = "  -0.02%  "

# Row 34
$ws.Range("D34").Value = "26.10"
$ws.Range("E34").Value = "  -7.04%  "

# Row 35
$ws.Range("E35").Value = "  -4.46%  "

# Row 36
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  -7.82%  "

# Row 37
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "53.60"
$ws.Range("E37").Value = "  -3.68%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0734"
$ws.Range("E38").Value = "  -6.71%  "

# Row 39
$ws.Range("D39").Value = "462.05"
$ws.Range("E39").Value = "  -6.46%  "

# Row 40
$ws.Range("E40").Value = "  -13.26%  "

# Row 41
$ws.Range("D41").Value = "0.0394"
$ws.Range("E41").Value = "  -8.03%  "

# Row 42
$ws.Range("E42").Value = "  -8.22%  "

# Row 43
$ws.Range("D43").Value = "8.41"
$ws.Range("E43").Value = "  -4.83%  "

# Row 44
$ws.Range("D44").Value = "2.841.23"
$ws.Range("E44").Value = "  -4.71%  "

# Row 45
$ws.Range("D45").Value = "0.265"
$ws.Range("E45").Value = "  -9.63%  "

# Row 46
$ws.Range("D46").Value = "2.26"
$ws.Range("E46").Value = "  -10.07%  "

# Row 47
$ws.Range("E47").Value = "  -9.01%  "

# Row 48
$ws.Range("E48").Value = "  -0.06%  "

# Row 49
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "2.31"
$ws.Range("E49").Value = "  -8.00%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.114"
$ws.Range("E50").Value = "  -4.82%  "

# Row 51
$ws.Range("D51").Value = "119.06"
